$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "TestCases" (sheet1): add a new row 3 (TC2) that mirrors row 2
# ------------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("TestCases")

[void]$wsTestCases.Range("A2:D2").Copy()
[void]$wsTestCases.Range("A3:D3").PasteSpecial(-4122)
$wsTestCases.Range("A3").Value = "TC2"
$wsTestCases.Range("B3").Value = "Y"
$wsTestCases.Range("C3").Value = "TC1_MercurySite|TC1_MercurySite_LoginTest"
$wsTestCases.Range("D3").Value = "TC1_MercurySite|CloseBrowser"
$wsTestCases.Rows("3:3").RowHeight = $wsTestCases.Rows("2:2").RowHeight

[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("C13").Select()

# ------------------------------------------------------------------
# Sheet "TestData" (sheet2): collapse the duplicated rows down to a
# TC1 / TC2 pair and turn TC2's username into a mailto hyperlink.
# ------------------------------------------------------------------
$wsTestData = $wb.Worksheets.Item("TestData")

[void]$wsTestData.Rows("3:6").Delete()

[void]$wsTestData.Range("A2:C2").Copy()
[void]$wsTestData.Range("A3:C3").PasteSpecial(-4122)

$wsTestData.Range("A3").Value = "TC2"
[void]$wsTestData.Hyperlinks.Add($wsTestData.Range("B3"), "mailto:santosh.pandhare@gmail.com", "", "", "santosh.pandhare@gmail.com")
$wsTestData.Range("C3").Value = "mercury"

$wsTestData.Columns("B:B").ColumnWidth = 27.59

[void]$wsTestData.Activate()
[void]$wsTestData.Range("I6").Select()

[void]$wb.Save()
